$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update published StructureDefinition metadata ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/evaluated-output"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet: keep in sync with the metadata URL, and drop the ---
# --- now-redundant root-level constraint text on the Extension row.   ---
$elem = $wb.Worksheets.Item("Elements")
$elem.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/evaluated-output"
$elem.Range("AI2").Value = ""
